# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously-last row (24) loses the "last row" date format and takes the
# regular date format used by the rest of the column.
$ws.Range("A24").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new day's data as row 25.
$ws.Range("A25").Value = 45765
$ws.Range("B25").Value = 98
$ws.Range("C25").Value = 101
$ws.Range("D25").Value = 99

# New last row takes on the distinct "last row" date number format.
$ws.Range("A25").NumberFormat = "YYYY-MM-DD"
